# Apply the three text edits described by the diff:
#   Slide 3 - "Rectangle 21" shape: EVOLUTION -> EVALUATION
#   Slide 4 - Title placeholder:   Gathering of Data -> Collection of Data
#   Slide 5 - Title placeholder:   Data preparation -> Data Preparation (split run)

$p = $ppt.ActivePresentation

# --- Slide 3: "EVOLUTION" -> "EVALUATION" -----------------------------
$slide3 = $p.Slides.Item(3)
for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
    $shape = $slide3.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "EVOLUTION") {
            $shape.TextFrame.TextRange.Text = "EVALUATION"
        }
    }
}

# --- Slide 4: "Gathering of Data" -> "Collection of Data" -------------
$slide4 = $p.Slides.Item(4)
$title4 = $slide4.Shapes.Item(1)
$title4.TextFrame.TextRange.Text = "Collection of Data"

# --- Slide 5: "Data preparation" -> "Data " + "Preparation" -----------
$slide5 = $p.Slides.Item(5)
$title5 = $slide5.Shapes.Item(1)
$title5.TextFrame.TextRange.Text = "Data "
$title5.TextFrame.TextRange.InsertAfter("Preparation") | Out-Null
